# Auto-generated Excel COM-interop script to apply price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 889.86664
$ws.Range("I8").Value = 395
$ws.Range("K8").Value = 1185
$ws.Range("M8").Value = -1046
$ws.Range("H47").Value = 24500
$ws.Range("I47").Value = 34000
$ws.Range("J47").Value = 15000
$ws.Range("K47").Value = 34000
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = -33028
$ws.Range("N47").Value = -16944
$ws.Range("H48").Value = 8650
$ws.Range("I48").Value = 475
$ws.Range("J48").Value = 25000
$ws.Range("K48").Value = 1425
$ws.Range("L48").Value = 75000
$ws.Range("M48").Value = -1133
$ws.Range("N48").Value = -75584
$ws.Range("H56").Value = 8650
$ws.Range("I56").Value = 475
$ws.Range("J56").Value = 25000
$ws.Range("K56").Value = 1425
$ws.Range("L56").Value = 75000
$ws.Range("M56").Value = -891
$ws.Range("N56").Value = -76068
$ws.Range("H131").Value = 5461.6665
$ws.Range("I131").Value = 2329.375
$ws.Range("K131").Value = 6988.125
$ws.Range("M131").Value = -1948.125
$ws.Range("H137").Value = 1113449.5
$ws.Range("I137").Value = 1668365.9
$ws.Range("K137").Value = 5005097.699999999
$ws.Range("M137").Value = -5002547.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4069
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 4069
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 4069
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -4295
$ws.Range("H74").Value = 2938.2354
$ws.Range("I74").Value = 1143.5217
$ws.Range("J74").Value = 6690.8184
$ws.Range("K74").Value = 1143.5217
$ws.Range("L74").Value = 6690.8184
$ws.Range("M74").Value = -269.5217
$ws.Range("N74").Value = -8438.8184
$ws.Range("H77").Value = 2938.2354
$ws.Range("I77").Value = 1143.5217
$ws.Range("J77").Value = 6690.8184
$ws.Range("K77").Value = 5717.6085
$ws.Range("L77").Value = 33454.092
$ws.Range("M77").Value = -1349.6085
$ws.Range("N77").Value = -42190.092
$ws.Range("H102").Value = 1854.4
$ws.Range("I102").Value = 1511.0555
$ws.Range("K102").Value = 1511.0555
$ws.Range("M102").Value = 110.9445000000001
$ws.Range("H116").Value = 4069
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4069
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4069
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -8657
$ws.Range("H122").Value = 1972.75
$ws.Range("I122").Value = 1866.6666
$ws.Range("K122").Value = 5599.9998
$ws.Range("M122").Value = -3149.9998
$ws.Range("H131").Value = 69000
$ws.Range("J131").Value = 69000
$ws.Range("L131").Value = 69000
$ws.Range("N131").Value = -79080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4069
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 4069
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 4069
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -4297
$ws.Range("H134").Value = 8166.5654
$ws.Range("I134").Value = 10021.814
$ws.Range("J134").Value = 5530.1577
$ws.Range("K134").Value = 30065.442
$ws.Range("L134").Value = 16590.4731
$ws.Range("M134").Value = -27530.442
$ws.Range("N134").Value = -21660.4731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2700.5386
$ws.Range("I31").Value = 1811.8125
$ws.Range("J31").Value = 4122.5
$ws.Range("K31").Value = 1811.8125
$ws.Range("L31").Value = 4122.5
$ws.Range("M31").Value = -1516.8125
$ws.Range("N31").Value = -4712.5
$ws.Range("H34").Value = 2700.5386
$ws.Range("I34").Value = 1811.8125
$ws.Range("J34").Value = 4122.5
$ws.Range("K34").Value = 1811.8125
$ws.Range("L34").Value = 4122.5
$ws.Range("M34").Value = -1609.8125
$ws.Range("N34").Value = -4526.5
$ws.Range("H86").Value = 4748
$ws.Range("I86").Value = 4997.3335
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 4997.3335
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -3874.3335
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 4748
$ws.Range("I89").Value = 4997.3335
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 24986.6675
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -19370.6675
$ws.Range("N89").Value = -31232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6126.857
$ws.Range("I5").Value = 6222.25
$ws.Range("K5").Value = 18666.75
$ws.Range("M5").Value = -18554.75
$ws.Range("H16").Value = 10750.25
$ws.Range("J16").Value = 17000.5
$ws.Range("L16").Value = 51001.5
$ws.Range("N16").Value = -51347.5
$ws.Range("H21").Value = 4076
$ws.Range("I21").Value = 1791.5
$ws.Range("J21").Value = 5599
$ws.Range("K21").Value = 5374.5
$ws.Range("L21").Value = 16797
$ws.Range("M21").Value = -5201.5
$ws.Range("N21").Value = -17143
$ws.Range("H37").Value = 139479.86
$ws.Range("J37").Value = 139479.86
$ws.Range("L37").Value = 418439.58
$ws.Range("N37").Value = -418663.58
$ws.Range("H68").Value = 2147.75
$ws.Range("I68").Value = 898
$ws.Range("J68").Value = 4230.6665
$ws.Range("K68").Value = 2694
$ws.Range("L68").Value = 12691.9995
$ws.Range("M68").Value = -1883
$ws.Range("N68").Value = -14313.9995
$ws.Range("H71").Value = 2147.75
$ws.Range("I71").Value = 898
$ws.Range("J71").Value = 4230.6665
$ws.Range("K71").Value = 8082
$ws.Range("L71").Value = 38075.9985
$ws.Range("M71").Value = -4026
$ws.Range("N71").Value = -46187.9985
$ws.Range("H114").Value = 2787
$ws.Range("I114").Value = 845
$ws.Range("J114").Value = 3272.5
$ws.Range("K114").Value = 2535
$ws.Range("L114").Value = 9817.5
$ws.Range("M114").Value = 719
$ws.Range("N114").Value = -16325.5
$ws.Range("H124").Value = 24852
$ws.Range("I124").Value = 23450
$ws.Range("K124").Value = 70350
$ws.Range("M124").Value = -65440
$ws.Range("H135").Value = 6126.857
$ws.Range("I135").Value = 6222.25
$ws.Range("K135").Value = 56000.25
$ws.Range("M135").Value = -53465.25
$ws.Range("H136").Value = 18084.047
$ws.Range("I136").Value = 9165.5
$ws.Range("J136").Value = 29975.445
$ws.Range("K136").Value = 27496.5
$ws.Range("L136").Value = 89926.33499999999
$ws.Range("M136").Value = -22396.5
$ws.Range("N136").Value = -100126.335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 18000
$ws.Range("I43").Value = 18000
$ws.Range("K43").Value = 18000
$ws.Range("M43").Value = -17849
$ws.Range("H57").Value = 27150
$ws.Range("I57").Value = 2250
$ws.Range("K57").Value = 2250
$ws.Range("M57").Value = -1430
$ws.Range("H80").Value = 2799
$ws.Range("I80").Value = 2799
$ws.Range("K80").Value = 2799
$ws.Range("M80").Value = -1801
$ws.Range("H83").Value = 2799
$ws.Range("I83").Value = 2799
$ws.Range("K83").Value = 13995
$ws.Range("M83").Value = -9003
$ws.Range("H122").Value = 3668.3333
$ws.Range("I122").Value = 3668.3333
$ws.Range("K122").Value = 11004.9999
$ws.Range("M122").Value = -8554.999899999999
$ws.Range("H136").Value = 43166.5
$ws.Range("J136").Value = 43166.5
$ws.Range("L136").Value = 129499.5
$ws.Range("N136").Value = -134599.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1349
$ws.Range("I22").Value = 588.2
$ws.Range("J22").Value = 2300
$ws.Range("K22").Value = 588.2
$ws.Range("L22").Value = 2300
$ws.Range("M22").Value = -293.2
$ws.Range("N22").Value = -2890
$ws.Range("H27").Value = 1349
$ws.Range("I27").Value = 588.2
$ws.Range("J27").Value = 2300
$ws.Range("K27").Value = 588.2
$ws.Range("L27").Value = 2300
$ws.Range("M27").Value = -481.2
$ws.Range("N27").Value = -2514
$ws.Range("H45").Value = 28999.857
$ws.Range("I45").Value = 25500
$ws.Range("K45").Value = 25500
$ws.Range("M45").Value = -25093
$ws.Range("H46").Value = 4763265.5
$ws.Range("I46").Value = 7693105
$ws.Range("J46").Value = 2276
$ws.Range("K46").Value = 7693105
$ws.Range("L46").Value = 2276
$ws.Range("M46").Value = -7692917
$ws.Range("N46").Value = -2652
$ws.Range("H68").Value = 2068.5264
$ws.Range("I68").Value = 2124.2942
$ws.Range("J68").Value = 1594.5
$ws.Range("K68").Value = 2124.2942
$ws.Range("L68").Value = 1594.5
$ws.Range("M68").Value = -1375.2942
$ws.Range("N68").Value = -3092.5
$ws.Range("H71").Value = 2068.5264
$ws.Range("I71").Value = 2124.2942
$ws.Range("J71").Value = 1594.5
$ws.Range("K71").Value = 10621.471
$ws.Range("L71").Value = 7972.5
$ws.Range("M71").Value = -6877.471
$ws.Range("N71").Value = -15460.5
$ws.Range("H122").Value = 5681.846
$ws.Range("I122").Value = 3233.1667
$ws.Range("K122").Value = 9699.500100000001
$ws.Range("M122").Value = -7249.500100000001
$ws.Range("H132").Value = 3284.6072
$ws.Range("I132").Value = 2788.9167
$ws.Range("K132").Value = 8366.750100000001
$ws.Range("M132").Value = -5836.750100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 40000
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 40000
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H96").Value = 3624
$ws.Range("I96").Value = 1999.5
$ws.Range("J96").Value = 4165.5
$ws.Range("K96").Value = 1999.5
$ws.Range("L96").Value = 4165.5
$ws.Range("M96").Value = -626.5
$ws.Range("N96").Value = -6911.5
$ws.Range("H107").Value = 9971.817999999999
$ws.Range("I107").Value = 17514.666
$ws.Range("J107").Value = 920.4
$ws.Range("K107").Value = 52543.99800000001
$ws.Range("L107").Value = 2761.2
$ws.Range("M107").Value = -50623.99800000001
$ws.Range("N107").Value = -6601.2
$ws.Range("H132").Value = 2031.4242
$ws.Range("I132").Value = 1811.2759
$ws.Range("K132").Value = 5433.8277
$ws.Range("M132").Value = -2903.8277
